$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(33, 8).Value = 203.66667  # H33
$ws.Cells.Item(33, 9).Value = 212.11765  # I33
$ws.Cells.Item(33, 11).Value = 212.11765  # K33
$ws.Cells.Item(33, 13).Value = 16.88235  # M33

$ws.Cells.Item(40, 8).Value = 1322.7  # H40
$ws.Cells.Item(40, 9).Value = 775  # I40
$ws.Cells.Item(40, 10).Value = 2144.25  # J40
$ws.Cells.Item(40, 11).Value = 775  # K40
$ws.Cells.Item(40, 12).Value = 2144.25  # L40
$ws.Cells.Item(40, 13).Value = -600  # M40
$ws.Cells.Item(40, 14).Value = -2494.25  # N40

$ws.Cells.Item(49, 8).Value = 573.75  # H49
$ws.Cells.Item(49, 10).Value = 800  # J49
$ws.Cells.Item(49, 12).Value = 2400  # L49
$ws.Cells.Item(49, 14).Value = -2672  # N49

$ws.Cells.Item(64, 8).Value = 3380  # H64
$ws.Cells.Item(64, 9).Value = 0  # I64
$ws.Cells.Item(64, 10).Value = 3380  # J64
$ws.Cells.Item(64, 11).Value = 0  # K64
$ws.Cells.Item(64, 12).Value = 3380  # L64
$ws.Cells.Item(64, 13).ClearContents()  # M64
$ws.Cells.Item(64, 14).Value = -3876  # N64

$ws.Cells.Item(67, 8).Value = 3380  # H67
$ws.Cells.Item(67, 9).Value = 0  # I67
$ws.Cells.Item(67, 10).Value = 3380  # J67
$ws.Cells.Item(67, 11).Value = 0  # K67
$ws.Cells.Item(67, 12).Value = 3380  # L67
$ws.Cells.Item(67, 13).ClearContents()  # M67
$ws.Cells.Item(67, 14).Value = -5096  # N67

$ws.Cells.Item(107, 8).Value = 720.03125  # H107
$ws.Cells.Item(107, 9).Value = 685.84  # I107
$ws.Cells.Item(107, 11).Value = 685.84  # K107
$ws.Cells.Item(107, 13).Value = 1234.16  # M107

$ws.Cells.Item(112, 8).Value = 4808730  # H112
$ws.Cells.Item(112, 10).Value = 1090.8334  # J112
$ws.Cells.Item(112, 12).Value = 3272.5002  # L112
$ws.Cells.Item(112, 14).Value = -5488.5002  # N112

$ws.Cells.Item(129, 8).Value = 753.3019  # H129
$ws.Cells.Item(129, 10).Value = 800.8542  # J129
$ws.Cells.Item(129, 12).Value = 2402.5626  # L129
$ws.Cells.Item(129, 14).Value = -12402.5626  # N129

$ws.Cells.Item(134, 8).Value = 47974.5  # H134
$ws.Cells.Item(134, 10).Value = 47974.5  # J134
$ws.Cells.Item(134, 12).Value = 47974.5  # L134
$ws.Cells.Item(134, 14).Value = -58114.5  # N134

$ws.Cells.Item(138, 8).Value = 2065.0266  # H138
$ws.Cells.Item(138, 9).Value = 1141.9286  # I138
$ws.Cells.Item(138, 10).Value = 2614.9575  # J138
$ws.Cells.Item(138, 11).Value = 3425.7858  # K138
$ws.Cells.Item(138, 12).Value = 7844.872499999999  # L138
$ws.Cells.Item(138, 13).Value = 1714.2142  # M138
$ws.Cells.Item(138, 14).Value = -18124.8725  # N138

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 5722.73  # H32
$ws.Cells.Item(32, 9).Value = 4604.094  # I32
$ws.Cells.Item(32, 10).Value = 12061.667  # J32
$ws.Cells.Item(32, 11).Value = 4604.094  # K32
$ws.Cells.Item(32, 12).Value = 12061.667  # L32
$ws.Cells.Item(32, 13).Value = -4317.094  # M32
$ws.Cells.Item(32, 14).Value = -12635.667  # N32

$ws.Cells.Item(45, 8).Value = 2672.818  # H45
$ws.Cells.Item(45, 9).Value = 2521  # I45
$ws.Cells.Item(45, 10).Value = 2938.5  # J45
$ws.Cells.Item(45, 11).Value = 2521  # K45
$ws.Cells.Item(45, 12).Value = 2938.5  # L45
$ws.Cells.Item(45, 13).Value = -2144  # M45
$ws.Cells.Item(45, 14).Value = -3692.5  # N45

$ws.Cells.Item(61, 8).Value = 1919.5135  # H61
$ws.Cells.Item(61, 9).Value = 1584.1333  # I61
$ws.Cells.Item(61, 11).Value = 1584.1333  # K61
$ws.Cells.Item(61, 13).Value = -1372.1333  # M61

$ws.Cells.Item(136, 8).Value = 1919.5135  # H136
$ws.Cells.Item(136, 9).Value = 1584.1333  # I136
$ws.Cells.Item(136, 11).Value = 4752.3999  # K136
$ws.Cells.Item(136, 13).Value = -2202.3999  # M136

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(86, 8).Value = 1882.1111  # H86
$ws.Cells.Item(86, 9).Value = 1759.125  # I86
$ws.Cells.Item(86, 10).Value = 2128.0833  # J86
$ws.Cells.Item(86, 11).Value = 1759.125  # K86
$ws.Cells.Item(86, 12).Value = 2128.0833  # L86
$ws.Cells.Item(86, 13).Value = -636.125  # M86
$ws.Cells.Item(86, 14).Value = -4374.0833  # N86

$ws.Cells.Item(89, 8).Value = 1882.1111  # H89
$ws.Cells.Item(89, 9).Value = 1759.125  # I89
$ws.Cells.Item(89, 10).Value = 2128.0833  # J89
$ws.Cells.Item(89, 11).Value = 8795.625  # K89
$ws.Cells.Item(89, 12).Value = 10640.4165  # L89
$ws.Cells.Item(89, 13).Value = -3179.625  # M89
$ws.Cells.Item(89, 14).Value = -21872.4165  # N89

$ws.Cells.Item(99, 8).Value = 1378.6842  # H99
$ws.Cells.Item(99, 9).Value = 1279.5  # I99
$ws.Cells.Item(99, 10).Value = 1488.8889  # J99
$ws.Cells.Item(99, 11).Value = 1279.5  # K99
$ws.Cells.Item(99, 12).Value = 1488.8889  # L99
$ws.Cells.Item(99, 13).Value = 218.5  # M99
$ws.Cells.Item(99, 14).Value = -4484.8889  # N99

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(99, 8).Value = 4038.5454  # H99
$ws.Cells.Item(99, 9).Value = 2871.75  # I99
$ws.Cells.Item(99, 11).Value = 2871.75  # K99
$ws.Cells.Item(99, 13).Value = -1373.75  # M99

$ws.Cells.Item(105, 8).Value = 733.5714  # H105
$ws.Cells.Item(105, 9).Value = 733.5714  # I105
$ws.Cells.Item(105, 11).Value = 733.5714  # K105
$ws.Cells.Item(105, 13).Value = 1013.4286  # M105

$ws.Cells.Item(107, 8).Value = 1017.26666  # H107
$ws.Cells.Item(107, 9).Value = 395.94446  # I107
$ws.Cells.Item(107, 10).Value = 1949.25  # J107
$ws.Cells.Item(107, 11).Value = 395.94446  # K107
$ws.Cells.Item(107, 12).Value = 1949.25  # L107
$ws.Cells.Item(107, 13).Value = 1524.05554  # M107
$ws.Cells.Item(107, 14).Value = -5789.25  # N107

$ws.Cells.Item(126, 8).Value = 4038.5454  # H126
$ws.Cells.Item(126, 9).Value = 2871.75  # I126
$ws.Cells.Item(126, 11).Value = 8615.25  # K126
$ws.Cells.Item(126, 13).Value = -6145.25  # M126

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 1171.2632  # H5
$ws.Cells.Item(5, 9).Value = 867.7308  # I5
$ws.Cells.Item(5, 10).Value = 1828.9166  # J5
$ws.Cells.Item(5, 11).Value = 2603.1924  # K5
$ws.Cells.Item(5, 12).Value = 5486.7498  # L5
$ws.Cells.Item(5, 13).Value = -2491.1924  # M5
$ws.Cells.Item(5, 14).Value = -5710.7498  # N5

$ws.Cells.Item(110, 8).Value = 1957.5  # H110
$ws.Cells.Item(110, 9).Value = 1957.5  # I110
$ws.Cells.Item(110, 11).Value = 5872.5  # K110
$ws.Cells.Item(110, 13).Value = -1782.5  # M110

$ws.Cells.Item(122, 8).Value = 897.63635  # H122
$ws.Cells.Item(122, 9).Value = 586.6667  # I122
$ws.Cells.Item(122, 10).Value = 946.7368  # J122
$ws.Cells.Item(122, 11).Value = 5280.0003  # K122
$ws.Cells.Item(122, 12).Value = 8520.6312  # L122
$ws.Cells.Item(122, 13).Value = -2830.0003  # M122
$ws.Cells.Item(122, 14).Value = -13420.6312  # N122

$ws.Cells.Item(131, 8).Value = 732.98  # H131
$ws.Cells.Item(131, 10).Value = 746.0928  # J131
$ws.Cells.Item(131, 12).Value = 2238.2784  # L131
$ws.Cells.Item(131, 14).Value = -12318.2784  # N131

$ws.Cells.Item(135, 8).Value = 1171.2632  # H135
$ws.Cells.Item(135, 9).Value = 867.7308  # I135
$ws.Cells.Item(135, 10).Value = 1828.9166  # J135
$ws.Cells.Item(135, 11).Value = 7809.577200000001  # K135
$ws.Cells.Item(135, 12).Value = 16460.2494  # L135
$ws.Cells.Item(135, 13).Value = -5274.577200000001  # M135
$ws.Cells.Item(135, 14).Value = -21530.2494  # N135

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(80, 8).Value = 3365.2173  # H80
$ws.Cells.Item(80, 10).Value = 3652.3333  # J80
$ws.Cells.Item(80, 12).Value = 3652.3333  # L80
$ws.Cells.Item(80, 14).Value = -5648.3333  # N80

$ws.Cells.Item(83, 8).Value = 3365.2173  # H83
$ws.Cells.Item(83, 10).Value = 3652.3333  # J83
$ws.Cells.Item(83, 12).Value = 18261.6665  # L83
$ws.Cells.Item(83, 14).Value = -28245.6665  # N83

$ws.Cells.Item(122, 8).Value = 3972.318  # H122
$ws.Cells.Item(122, 9).Value = 4026.6365  # I122
$ws.Cells.Item(122, 10).Value = 3918  # J122
$ws.Cells.Item(122, 11).Value = 12079.9095  # K122
$ws.Cells.Item(122, 12).Value = 11754  # L122
$ws.Cells.Item(122, 13).Value = -9629.9095  # M122
$ws.Cells.Item(122, 14).Value = -16654  # N122

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 4431.3706  # H7
$ws.Cells.Item(7, 9).Value = 4307.35  # I7
$ws.Cells.Item(7, 11).Value = 4307.35  # K7
$ws.Cells.Item(7, 13).Value = -4195.35  # M7

$ws.Cells.Item(40, 8).Value = 4043.9412  # H40
$ws.Cells.Item(40, 9).Value = 4017.6428  # I40
$ws.Cells.Item(40, 10).Value = 4166.6665  # J40
$ws.Cells.Item(40, 11).Value = 4017.6428  # K40
$ws.Cells.Item(40, 12).Value = 4166.6665  # L40
$ws.Cells.Item(40, 13).Value = -3881.6428  # M40
$ws.Cells.Item(40, 14).Value = -4438.6665  # N40

$ws.Cells.Item(55, 8).Value = 165.5  # H55
$ws.Cells.Item(55, 9).Value = 146  # I55
$ws.Cells.Item(55, 10).Value = 198  # J55
$ws.Cells.Item(55, 11).Value = 146  # K55
$ws.Cells.Item(55, 12).Value = 198  # L55
$ws.Cells.Item(55, 13).Value = 27  # M55
$ws.Cells.Item(55, 14).Value = -544  # N55

$ws.Cells.Item(122, 8).Value = 2183091.5  # H122
$ws.Cells.Item(122, 9).Value = 2455015.5  # I122
$ws.Cells.Item(122, 10).Value = 7700  # J122
$ws.Cells.Item(122, 11).Value = 7365046.5  # K122
$ws.Cells.Item(122, 12).Value = 23100  # L122
$ws.Cells.Item(122, 13).Value = -7362596.5  # M122
$ws.Cells.Item(122, 14).Value = -28000  # N122

$ws.Cells.Item(126, 8).Value = 4431.3706  # H126
$ws.Cells.Item(126, 9).Value = 4307.35  # I126
$ws.Cells.Item(126, 11).Value = 12922.05  # K126
$ws.Cells.Item(126, 13).Value = -10452.05  # M126

$ws.Cells.Item(132, 8).Value = 379091.25  # H132
$ws.Cells.Item(132, 9).Value = 465274.7  # I132
$ws.Cells.Item(132, 10).Value = 5629.6665  # J132
$ws.Cells.Item(132, 11).Value = 1395824.1  # K132
$ws.Cells.Item(132, 12).Value = 16888.9995  # L132
$ws.Cells.Item(132, 13).Value = -1393294.1  # M132
$ws.Cells.Item(132, 14).Value = -21948.9995  # N132

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(62, 8).Value = 4178  # H62
$ws.Cells.Item(62, 10).Value = 4700.3335  # J62
$ws.Cells.Item(62, 12).Value = 4700.3335  # L62
$ws.Cells.Item(62, 14).Value = -5948.3335  # N62

$ws.Cells.Item(65, 8).Value = 4178  # H65
$ws.Cells.Item(65, 10).Value = 4700.3335  # J65
$ws.Cells.Item(65, 12).Value = 23501.6675  # L65
$ws.Cells.Item(65, 14).Value = -29741.6675  # N65

$ws.Cells.Item(81, 8).Value = 1592.2727  # H81
$ws.Cells.Item(81, 9).Value = 427  # I81
$ws.Cells.Item(81, 10).Value = 3631.5  # J81
$ws.Cells.Item(81, 11).Value = 854  # K81
$ws.Cells.Item(81, 12).Value = 7263  # L81
$ws.Cells.Item(81, 13).Value = 207  # M81
$ws.Cells.Item(81, 14).Value = -9385  # N81

$ws.Cells.Item(84, 8).Value = 1592.2727  # H84
$ws.Cells.Item(84, 9).Value = 427  # I84
$ws.Cells.Item(84, 10).Value = 3631.5  # J84
$ws.Cells.Item(84, 11).Value = 4270  # K84
$ws.Cells.Item(84, 12).Value = 36315  # L84
$ws.Cells.Item(84, 13).Value = 1034  # M84
$ws.Cells.Item(84, 14).Value = -46923  # N84

$ws.Cells.Item(132, 8).Value = 1317.2  # H132
$ws.Cells.Item(132, 9).Value = 1031.1666  # I132
$ws.Cells.Item(132, 11).Value = 3093.4998  # K132
$ws.Cells.Item(132, 13).Value = -563.4998000000001  # M132
